# Regenerate s_vals data to filter save games: update numeric cell values
# in B2:E7 and G2:G7 (F column / "Win" flags remain unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 3.223369029078222;  E = 0.5333859586016987;  G = 8.656069925401464 }
    3 = @{ B = 1.445647641019636;  C = 0.04103571897497393; D = 3.223369029078222;  E = 13.86384647080068;   G = 18.57389885987352 }
    4 = @{ B = 0.2881169905109251; C = 0.3048912486333797;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 1.84748871573303 }
    5 = @{ B = 1.445647641019636;  C = 1.626987699542094;   D = 0.1496068669990043; E = 0.5333859586016987;  G = 3.755628166162433 }
    6 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 18.71679738969934;  E = 0.5333859586016987;  G = 24.14949828602258 }
    7 = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.1496068669990043; E = 0.5333859586016987;  G = 5.582307763322248 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("B$row").Value = $rowVals.B
    $ws.Range("C$row").Value = $rowVals.C
    $ws.Range("D$row").Value = $rowVals.D
    $ws.Range("E$row").Value = $rowVals.E
    $ws.Range("G$row").Value = $rowVals.G
}
